$d = $word.ActiveDocument

# Bump the FreeBSD / XigmaNAS release version referenced throughout the
# manual from 12.1 to 12.2 (title, svn checkout URL, cd target, the three
# "cp .../xigmanas/12.1/..." loader copy commands, the bookmarked cp
# command, and the closing note).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Text = "12.1"
$find.Replacement.Text = "12.2"
$find.Forward = $true
$find.Wrap = 1                 # wdFindContinue
$find.Format = $false
$find.MatchCase = $false
$find.MatchWholeWord = $false
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute("12.1", $false, $false, $false, $false, $false, $true, 1, $false, "12.2", 2) | Out-Null
